$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp update in the header cell
$ws.Range("A1").Value = "Datos actualizados a 7 de Mayo de 2020 a las 05:03"

# Mexico row (23) - updated recovered/active counts
$ws.Range("D23").Value = 17781
$ws.Range("E23").Value = 7149

# Row 25: now Pakistan (was Suecia) with new stats
$ws.Range("A25").Value = "Pakistan"
$ws.Range("B25").Value = 24073
$ws.Range("C25").Value = 859
$ws.Range("D25").Value = 6464
$ws.Range("E25").Value = 17045
$ws.Range("F25").Value = 111
$ws.Range("G25").Value = 20
$ws.Range("H25").Value = 564

# Row 26: now Suecia (was Pakistan) inherits previous Suecia stats
$ws.Range("A26").Value = "Suecia"
$ws.Range("B26").Value = 23918
$ws.Range("D26").Value = 4074
$ws.Range("E26").Value = 16903
$ws.Range("F26").Value = 425
$ws.Range("H26").Value = 2941

# Row 77: now Bolivia (was Guinea) with new stats
$ws.Range("A77").Value = "Bolivia"
$ws.Range("B77").Value = 1886
$ws.Range("C77").Value = 84
$ws.Range("D77").Value = 198
$ws.Range("E77").Value = 1597
$ws.Range("F77").Value = 3
$ws.Range("G77").Value = 5
$ws.Range("H77").Value = 91

# Row 78: now Guinea (was Bolivia) inherits previous Guinea stats
$ws.Range("A78").Value = "Guinea"
$ws.Range("B78").Value = 1856
$ws.Range("D78").Value = 597
$ws.Range("E78").Value = 1248
$ws.Range("F78").Value = 0
$ws.Range("H78").Value = 11

# Row 171: minor data correction
$ws.Range("D171").Value = 55
$ws.Range("E171").Value = 5

# Row 191: now Nueva Caledonia (was Belice) with new stats
$ws.Range("A191").Value = "Nueva Caledonia"
$ws.Range("D191").Value = 18
$ws.Range("F191").Value = 1
$ws.Range("H191").Value = 0

# Row 192: now Belice (was Nueva Caledonia) inherits previous Belice stats
$ws.Range("A192").Value = "Belice"
$ws.Range("D192").Value = 16
$ws.Range("F192").Value = 0
$ws.Range("H192").Value = 2

# Row 198: now Dominica (was Curazao) with new stats
$ws.Range("A198").Value = "Dominica"
$ws.Range("D198").Value = 14
$ws.Range("H198").Value = 0

# Row 199: now Curazao (was Dominica) inherits previous Curazao stats
$ws.Range("A199").Value = "Curazao"
$ws.Range("D199").Value = 13
$ws.Range("H199").Value = 1
